$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the title
#    (first) paragraph: empty run + bold "Meta description" run +
#    normal run with the rest of the sentence.
# -----------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter()
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Style = "Normal"

$metaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Play the Bier Haus Oktoberfest online slot game for free and discover its multiple bonuses and potential big payouts. Review of the game features and graphics.</w:t></w:r></w:p>'
$metaFull = $d.Range($metaPara.Range.Start, $metaPara.Range.End)
$null = $metaFull.InsertXML($metaXml)

# -----------------------------------------------------------------
# 2) Remove the trailing bold "Play Bier Haus Oktoberfest Slot for
#    Free - Review and Features" paragraph near the end of the
#    document, and replace the text of the remaining (italic) meta
#    description paragraph with the DALLE image prompt, keeping the
#    italic formatting and leading empty run intact.
# -----------------------------------------------------------------
$count = $d.Paragraphs.Count
$boldPara = $d.Paragraphs.Item($count - 1)
$boldRange = $d.Range($boldPara.Range.Start, $boldPara.Range.End)
$boldRange.Delete()

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastRange = $d.Range($lastPara.Range.Start, $lastPara.Range.End)
$lastRange.Text = "Prompt for DALLE: Create a fun and eye-catching feature image for Bier Haus Oktoberfest online slot game. The image should be in a cartoon style and should feature a happy Maya warrior with glasses. The image should have a festive Oktoberfest theme, with beer mugs, musical instruments, and typical Bavarian hats. The Maya warrior should be holding a beer mug and surrounded by Heidi and Hans, the game symbols. The image should include the game logo and convey the fun and exciting nature of the game. Use bold, bright colors to grab the viewer's attention and make them want to try their luck with Bier Haus Oktoberfest."
